$d = $word.ActiveDocument

# Update the ID placeholder text and remove the trailing space run by
# matching the full original text (including the trailing space) and
# replacing it with the new text (no trailing space).
[void]$d.Content.Find.Execute("**ID__AFFARS_pgi_5347_topic_3__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5347_301_1__ID**", 2)

# Adjust the first paragraph's formatting: add a paragraph border (space-only,
# no visible line) and change the left indent from 120 twips (6pt) to 225
# twips (11.25pt).
$p = $d.Paragraphs(1)
$pf = $p.Range.ParagraphFormat

$borders = $pf.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

$pf.LeftIndent = 11.25
